# Edit script: adds two "asesoria" (advising session) records into the
# "asesorias" sheet, and updates the "Usuario" values for the teacher
# (Maryem Ruiz, docentes sheet) and mentor (Emanuel Valencia, mentores sheet)
# records so each user's advising sessions can be tied back to their own
# account.

$wb = $excel.ActiveWorkbook

$docentes   = $wb.Worksheets.Item("docentes")
$mentores   = $wb.Worksheets.Item("mentores")
$asesorias  = $wb.Worksheets.Item("asesorias")

# Differentiate the "Usuario" login for the teacher and the mentor so that
# advising sessions can be linked to the specific user that created them.
$docentes.Range("D2").Value = "ss"
$mentores.Range("D2").Value = "sss"

# Row 2: Daniel Henao (student) <-> Maryem Ruiz (teacher, usuario "ss")
$asesorias.Range("A2").Value = "Daniel Henao"
$asesorias.Range("B2").Value = "s"
$asesorias.Range("C2").Value = "ss"
$asesorias.Range("D2").Value = "Maryem Ruiz"
$asesorias.Range("E2").Value = "Consulta general"
# "11-11-2023" looks like a date, so it would normally be auto-converted to
# a date serial by Excel's smart entry; build it as a text formula instead,
# then convert that formula to its plain static (text) value in place.
$asesorias.Range("F2").Formula = "=""11-11-2023"""
$asesorias.Range("F2").Copy()
$asesorias.Range("F2").PasteSpecial(-4163)
$asesorias.Range("G2").Value = "17:20 - 17:40"

# Row 3: Daniel Henao (student) <-> Emanuel Valencia (mentor, usuario "sss")
$asesorias.Range("A3").Value = "Daniel Henao"
$asesorias.Range("B3").Value = "s"
$asesorias.Range("C3").Value = "sss"
$asesorias.Range("D3").Value = "Emanuel Valencia"
$asesorias.Range("E3").Value = "Consulta general"
$asesorias.Range("F3").Formula = "=""11-11-2023"""
$asesorias.Range("F3").Copy()
$asesorias.Range("F3").PasteSpecial(-4163)
$asesorias.Range("G3").Value = "17:20 - 17:40"

$excel.CutCopyMode = $false

# Restore the active selection on the asesorias sheet similarly to the
# authored workbook.
$asesorias.Range("C4").Select()
